$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.307.24'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.17%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.865.54'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.12%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '234.65'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -0.84%  '
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +0.09%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4694'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.24%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2856'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -1.94%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06566'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +0.19%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '21.31'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -2.55%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07837'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -1.16%  '
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -1.05%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.881.32'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +0.90%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6978'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +2.39%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.087'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -1.46%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '268.97'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +0.80%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '30.272.77'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +0.06%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '13.78'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +0.00%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007624'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +2.64%  '
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +0.04%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.110.13'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.19%  '
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +0.02%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.227'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -1.64%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.150'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -0.47%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.432'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +2.13%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '167.08'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -0.22%  '
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -0.18%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.940'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -0.86%  '
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -1.78%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.09912'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +0.58%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.351'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -0.55%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.457'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -0.87%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.044'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -0.26%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04725'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +0.23%  '
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +0.22%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7021'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -0.22%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.718'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +0.53%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01870'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -0.52%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.743'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +5.01%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.326'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +0.61%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '72.69'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -2.08%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.948'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -0.19%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.4167'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +0.08%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.8346'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -1.37%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '103.17'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -0.21%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '969.86'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +1.81%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.097'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -1.08%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.142'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -1.18%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '34.43'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +0.82%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05682'
